$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (task submitted by Ahmed Omar)
$ws.Range("A2").Value = "Ahmed Omar Zakryia Albanna"
$ws.Range("B2").Value = "ahmedalbanna685@gmail.com"
$ws.Range("C2").Value = "https://github.com/aahmedd38/Security-Task.git"

# Turn the email and repo link cells into real hyperlinks
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:ahmedalbanna685@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/aahmedd38/Security-Task.git") | Out-Null

# Restore the selection that ends up being saved with the workbook
$ws.Range("C5").Select() | Out-Null
